$d = $word.ActiveDocument

# 1. "esse registro deve" -> "esse registro, além de"
$d.Content.Find.Execute("esse registro deve", $true, $false, $false, $false, $false,
                         $true, 1, $false, "esse registro, além de", 2)

# 2. "Além disso, atividades de planejamento são muito diversificadas, para atender todos os detalhes igualmente, as companhias as dividem em vários setores."
#    -> "As atividades do planejamento são muito diversificadas, para atender igualmente todas as partes, as companhias as dividem em setores."
$d.Content.Find.Execute("Além disso, atividades de planejamento são muito diversificadas, para atender todos os detalhes igualmente, as companhias as dividem em vários setores.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "As atividades do planejamento são muito diversificadas, para atender igualmente todas as partes, as companhias as dividem em setores.", 2)

# 3. " até a próxima troca." -> " e já se programam para a próxima troca."
$d.Content.Find.Execute(" até a próxima troca.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " e já se programam para a próxima troca.", 2)
